# Updates odds values on the "Jogos da Semana" worksheet to reflect
# refreshed FlashScore data (commit: "Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "H2"  = 4.2
    "I2"  = 6.25
    "K2"  = 2.3
    "O2"  = 1.25
    "P2"  = 3.75
    "Q2"  = 1.83
    "R2"  = 2.03
    "S2"  = 1.36
    "T2"  = 3
    "U2"  = 1.83
    "V2"  = 1.83
    "W2"  = 7
    "X2"  = 7
    "AB2" = 26
    "AC2" = 12
    "AG2" = 301
    "AP2" = 19
    "AQ2" = 23
    "AT2" = 3
    "AX2" = 34
    "AZ2" = 126

    "Q3" = 2.5
    "R3" = 1.5

    "G4" = 4.33
    "M4" = 1.13
    "N4" = 6

    "Q5" = 2.4
    "R5" = 1.53

    "G6"  = 3.25
    "H6"  = 2.9
    "I6"  = 2.3
    "J6"  = 3.95
    "K6"  = 1.91
    "L6"  = 2.92
    "N6"  = 5.3
    "O6"  = 1.5
    "P6"  = 2.25
    "Q6"  = 2.45
    "T6"  = 2.22
    "U6"  = 2.05
    "W6"  = 7.1
    "X6"  = 15
    "Y6"  = 12.5
    "Z6"  = 45
    "AA6" = 40
    "AB6" = 60
    "AC6" = 6.2
    "AH6" = 5.9
    "AI6" = 9.75
    "AJ6" = 9.75
    "AK6" = 23
    "AL6" = 23
    "AN6" = 4.9
    "AO6" = 19.5
    "AP6" = 32
    "AQ6" = 110
    "AR6" = 175
    "AT6" = 2.2
    "AU6" = 7.8
    "AW6" = 3.95
    "AX6" = 12
    "AY6" = 24
    "AZ6" = 55
    "BA6" = 110

    "G7"  = 2.47
    "I7"  = 2.77
    "J7"  = 3.05
    "L7"  = 3.3
    "W7"  = 7.2
    "X7"  = 11.5
    "Y7"  = 9.75
    "Z7"  = 27
    "AA7" = 23
    "AB7" = 35
    "AF7" = 75
    "AH7" = 8.25
    "AI7" = 14
    "AL7" = 25
    "AM7" = 35
    "AP7" = 21
    "AQ7" = 60
    "AR7" = 90
    "AS7" = 300
    "AT7" = 2.47
    "AU7" = 6.9
    "AX7" = 14.5
    "AY7" = 22
    "AZ7" = 65
    "BA7" = 100
    "BB7" = 250
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
